$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Record"
$ws.Range("B20").Value = "RJ Record"
$ws.Range("C20").Value = "Governo"
$ws.Range("D20").Value = "2025-04-02T18:38"
$ws.Range("E20").Value = "Positivo"
$ws.Range("F20").Value = "Locutor fala sobre a entrega de várias obras em Campos pelos 190 anos de elevação de vila à categoria de cidade"

$ws.Range("A21").Value = "Record"
$ws.Range("B21").Value = "RJ Record"
$ws.Range("C21").Value = "CCZ"
$ws.Range("D21").Value = "2025-04-02T18:38"
$ws.Range("E21").Value = "Positivo"
$ws.Range("F21").Value = "Para evitar acidentes, CCZ apreende animais soltos em vias públicas. Repórter *ao vivo*. Depoimento do diretor do CCZ, Carlos Morales. "
